# Small bug fixes: nudge a handful of label textboxes (tx9..tx18) that sit
# inside the diagram group on slide 1. Only their position changes
# (a:off x/y) - sizes (a:ext) are untouched.
#
# Shape.Left/Shape.Top are expressed in points, while the underlying OOXML
# stores EMUs (1 pt = 12700 EMU). The PowerPoint object model stores
# Left/Top as single-precision floats, so the literals below are chosen
# (to 1e-6 pt granularity) so that after the float32 round-trip the
# resulting EMU value lands exactly on the target from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(2)

$moves = @(
    @{ Name = "tx9";  Left = 354.7448888897638;  Top = 193.74408748818897 },
    @{ Name = "tx10"; Left = 382.72804349606304; Top = 217.86881289763778 },
    @{ Name = "tx11"; Left = 413.0984959370079;  Top = 241.58377852755905 },
    @{ Name = "tx12"; Left = 459.58290138582674; Top = 265.70851193700787 },
    @{ Name = "tx13"; Left = 498.7379617559055;  Top = 283.6948091496063  },
    @{ Name = "tx14"; Left = 500.14668291338586; Top = 311.01213098425194 },
    @{ Name = "tx15"; Left = 395.85542307086615; Top = 424.9202432204724  },
    @{ Name = "tx16"; Left = 404.8622894645669;  Top = 452.23756505511807 },
    @{ Name = "tx17"; Left = 277.35557655118106; Top = 288.11045844094485 },
    @{ Name = "tx18"; Left = 284.4659882519685;  Top = 315.4277195354331  }
)

foreach ($m in $moves) {
    $sh = $grp.GroupItems.Item($m.Name)
    $sh.Left = $m.Left
    $sh.Top = $m.Top
}
